$d = $word.ActiveDocument

$pairs = @(
    @("463÷7=66, 1", "473÷4=118, 1"),
    @("708÷5=141, 3", "234÷3=78, 0"),
    @("414÷7=59, 1", "588÷7=84, 0"),
    @("173÷9=19, 2", "942÷5=188, 2"),
    @("324÷8=40, 4", "134÷9=14, 8"),
    @("429÷6=71, 3", "153÷3=51, 0"),
    @("526÷5=105, 1", "478÷4=119, 2"),
    @("455÷5=91, 0", "909÷6=151, 3"),
    @("437÷7=62, 3", "875÷9=97, 2"),
    @("489÷2=244, 1", "627÷5=125, 2"),
    @("883÷2=441, 1", "599÷3=199, 2"),
    @("340÷5=68, 0", "495÷8=61, 7"),
    @("580÷3=193, 1", "232÷7=33, 1"),
    @("310÷2=155, 0", "804÷5=160, 4"),
    @("647÷8=80, 7", "723÷9=80, 3"),
    @("519÷5=103, 4", "907÷6=151, 1"),
    @("123÷4=30, 3", "909÷9=101, 0"),
    @("499÷9=55, 4", "104÷5=20, 4"),
    @("548÷2=274, 0", "927÷6=154, 3"),
    @("833÷2=416, 1", "193÷9=21, 4"),
    @("787÷9=87, 4", "799÷6=133, 1"),
    @("894÷2=447, 0", "971÷6=161, 5"),
    @("173÷4=43, 1", "668÷4=167, 0"),
    @("519÷4=129, 3", "600÷8=75, 0"),
    @("761÷5=152, 1", "817÷5=163, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
